$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 26.499453
$arr[0,7] = 79.49835900000001
$arr[0,8] = 0.7877954840311897
$arr[0,9] = 0.7877954840311898
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 32.09065066666667
$arr[0,13] = 96.271952
$arr[0,14] = 0.946743898370603
$arr[0,15] = 0.9467438983706029
$arr[0,16] = 850.3846890807521
$arr[0,17] = 7653.462201726768
$arr[0,18] = 0.7458405676704447
$arr[0,19] = 0.7458405676704447
$ws.Range("A2:T2").Value = $arr

# Row 3
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 26.499453
$arr[0,7] = 79.49835900000001
$arr[0,8] = 0.7877954840311897
$arr[0,9] = 0.7877954840311898
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.154570666666667
$arr[0,13] = 3.463712
$arr[0,14] = 0.03406234249527876
$arr[0,15] = 0.03406234249527877
$arr[0,16] = 30.595491116512
$arr[0,17] = 275.3594200486081
$arr[0,18] = 0.0268341595933043
$arr[0,19] = 0.0268341595933043
$ws.Range("A3:T3").Value = $arr

# Row 4
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "sCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 26.499453
$arr[0,7] = 79.49835900000001
$arr[0,8] = 0.7877954840311897
$arr[0,9] = 0.7877954840311898
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.650588
$arr[0,13] = 1.951764
$arr[0,14] = 0.01919375913411833
$arr[0,15] = 0.01919375913411833
$arr[0,16] = 17.240226128364
$arr[0,17] = 155.162035155276
$arr[0,18] = 0.01512075676744082
$arr[0,19] = 0.01512075676744082
$ws.Range("A4:T4").Value = $arr

# Row 5
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.993807666666666
$arr[0,7] = 17.981423
$arr[0,8] = 0.1781883804149286
$arr[0,9] = 0.1781883804149287
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 32.09065066666667
$arr[0,13] = 96.271952
$arr[0,14] = 0.946743898370603
$arr[0,15] = 0.9467438983706029
$arr[0,16] = 192.3451879941884
$arr[0,17] = 1731.106691947696
$arr[0,18] = 0.1686987619183735
$arr[0,19] = 0.1686987619183735
$ws.Range("A5:T5").Value = $arr

# Row 6
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.993807666666666
$arr[0,7] = 17.981423
$arr[0,8] = 0.1781883804149286
$arr[0,9] = 0.1781883804149287
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.154570666666667
$arr[0,13] = 3.463712
$arr[0,14] = 0.03406234249527876
$arr[0,15] = 0.03406234249527877
$arr[0,16] = 6.92027451357511
$arr[0,17] = 62.282470622176
$arr[0,18] = 0.006069513642372321
$arr[0,19] = 0.006069513642372324
$ws.Range("A6:T6").Value = $arr

# Row 7
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "sCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.993807666666666
$arr[0,7] = 17.981423
$arr[0,8] = 0.1781883804149286
$arr[0,9] = 0.1781883804149287
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.650588
$arr[0,13] = 1.951764
$arr[0,14] = 0.01919375913411833
$arr[0,15] = 0.01919375913411833
$arr[0,16] = 3.899499342241334
$arr[0,17] = 35.095494080172
$arr[0,18] = 0.003420104854182789
$arr[0,19] = 0.003420104854182789
$ws.Range("A7:T7").Value = $arr

# Row 8
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "sCs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.144217
$arr[0,7] = 3.432650999999999
$arr[0,8] = 0.03401613555388164
$arr[0,9] = 0.03401613555388164
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 32.09065066666667
$arr[0,13] = 96.271952
$arr[0,14] = 0.946743898370603
$arr[0,15] = 0.9467438983706029
$arr[0,16] = 36.71866803386133
$arr[0,17] = 330.4680123047519
$arr[0,18] = 0.03220456878178478
$arr[0,19] = 0.03220456878178478
$ws.Range("A8:T8").Value = $arr

# Row 9
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "sCs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.144217
$arr[0,7] = 3.432650999999999
$arr[0,8] = 0.03401613555388164
$arr[0,9] = 0.03401613555388164
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.154570666666667
$arr[0,13] = 3.463712
$arr[0,14] = 0.03406234249527876
$arr[0,15] = 0.03406234249527877
$arr[0,16] = 1.321079384501333
$arr[0,17] = 11.889714460512
$arr[0,18] = 0.001158669259602146
$arr[0,19] = 0.001158669259602146
$ws.Range("A9:T9").Value = $arr

# Row 10
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "sCs"
$arr[0,1] = "Cd80"
$arr[0,2] = "Cd274"
$arr[0,3] = "sCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.144217
$arr[0,7] = 3.432650999999999
$arr[0,8] = 0.03401613555388164
$arr[0,9] = 0.03401613555388164
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.650588
$arr[0,13] = 1.951764
$arr[0,14] = 0.01919375913411833
$arr[0,15] = 0.01919375913411833
$arr[0,16] = 0.744413849596
$arr[0,17] = 6.699724646363999
$arr[0,18] = 0.0006528975124947231
$arr[0,19] = 0.0006528975124947231
$ws.Range("A10:T10").Value = $arr

Write-Output "done"